# CAP018_MaintainBooking_TestData.xlsx - "Change in excel infra"
# Replace every occurrence of the origin/destination code "ANC" with "SFO"
# in columns A and B (rows 2-60) of the CAP018_BKG_00001 sheet, and move
# the active selection to F9 (no more frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAP018_BKG_00001")

for ($r = 2; $r -le 60; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value2 -eq "ANC") {
        $cellA.Value = "SFO"
    }

    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value2 -eq "ANC") {
        $cellB.Value = "SFO"
    }
}

# Update the sheet's view: scroll back to the top and move the selection
# to F9 (previously topLeftCell A13 / selection C4).
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("F9").Select()
